# TimeSheet.xlsx update
# - Extend the last section's date-range header from a single month
#   ("Esfand 98") to a multi-month span ("Esfand 98 - Khordad 99").
# - Add a new task line ("Documents (Read/Write)") with its hours,
#   bump the "Code Refactoring" hours, and let the section's totals
#   recalculate.
# - Re-color the previous section's header / totals row to match the
#   newly-active (latest) section's highlight color.
# - Widen column A so the longer Farsi date range fits, and leave the
#   selection/scroll position where the user ended up working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section header text: single month -> month range ---
$ws.Range("A77").Value = "اسفند 98 تا خرداد 99"

# --- Hours for the existing "Code Refactoring" task ---
$ws.Range("C79").Value = 3

# --- New task row: "Documents (Read/Write)" ---
$ws.Range("B80").Value = "Documents (Read/Write)"
$ws.Range("C80").Value = 3

# C87 = SUM(C78:C86) recalculates automatically (2 -> 7)

# --- Manually tracked @Home hours total for this section ---
$ws.Range("D89").Value = 7

# --- Re-color the now-previous section's header/footer cells to the
#     "active section" highlight (same fill used by the new last section) ---
$ws.Range("A62").Interior.Color = 5296274
$ws.Range("B72").Interior.Color = 5296274
$ws.Range("C72").Interior.Color = 5296274

# --- Column A needs to be a bit wider for the longer header text ---
$ws.Columns("A").ColumnWidth = 20.83

# --- Leave the view where the user finished editing ---
$ws.Range("E91").Select()
